# Append 12 new months (Dec-15 .. Nov-16) of indirect cost data to the
# existing table, extending the sheet from rows 1:12 to rows 1:24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Dec-15", 1, 258650.74, 55717.72, 314368.46),
    @("Jan-16", 1, 252175.98, 58329.03, 310505.01),
    @("Feb-16", 1, 244305.23, 60012.48, 304317.71),
    @("Mar-16", 1, 241621.2,  58528.24, 300149.44),
    @("Apr-16", 1, 239032.88, 63694.36, 302727.24),
    @("May-16", 1, 239635.18, 57655.44, 297290.62),
    @("Jun-16", 1, 243293.68, 54036.96, 297330.64),
    @("Jul-16", 1, 242000.98, 50080.31, 292081.29),
    @("Aug-16", 1, 234265.6,  54302.69, 288568.29),
    @("Sep-16", 1, 240823.88, 56481.26, 297305.14),
    @("Oct-16", 1, 240985.37, 56440.75, 297426.12),
    @("Nov-16", 1, 235683.88, 55490.59, 291174.47)
)

$startRow = 13
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]

    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
}
